# Auto-generated edit script: updates currentAveragePrice / Leve profit
# columns (H-N) for specific rows across the 8 job sheets, per the
# scheduled-runner price refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 55
$ws.Range("H55").Value = 542.5
$ws.Range("I55").Value = 593
$ws.Range("J55").Value = 492
$ws.Range("K55").Value = 593
$ws.Range("L55").Value = 492
$ws.Range("M55").Value = -379
$ws.Range("N55").Value = -920

# Row 111
$ws.Range("H111").Value = 3221.75
$ws.Range("I111").Value = 994.5
$ws.Range("J111").Value = 5449
$ws.Range("K111").Value = 2983.5
$ws.Range("L111").Value = 16347
$ws.Range("M111").Value = 83.5
$ws.Range("N111").Value = -22481

# Row 132
$ws.Range("H132").Value = 1237.3334
$ws.Range("I132").Value = 1237.3334
$ws.Range("K132").Value = 3712.0002
$ws.Range("M132").Value = -1182.0002

# Row 133
$ws.Range("H133").Value = 73999
$ws.Range("J133").Value = 73999
$ws.Range("L133").Value = 73999
$ws.Range("N133").Value = -84119

# Row 137
$ws.Range("H137").Value = 2739.75
$ws.Range("I137").Value = 2248
$ws.Range("K137").Value = 6744
$ws.Range("M137").Value = -4194

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1254.4546
$ws.Range("I2").Value = 699.8889
$ws.Range("K2").Value = 699.8889
$ws.Range("M2").Value = -586.8889

# Row 61
$ws.Range("H61").Value = 5308.4375
$ws.Range("I61").Value = 1849.4286
$ws.Range("K61").Value = 1849.4286
$ws.Range("M61").Value = -1637.4286

# Row 74
$ws.Range("H74").Value = 2634.1428
$ws.Range("I74").Value = 1969.4
$ws.Range("J74").Value = 4296
$ws.Range("K74").Value = 1969.4
$ws.Range("L74").Value = 4296
$ws.Range("M74").Value = -1095.4
$ws.Range("N74").Value = -6044

# Row 77
$ws.Range("H77").Value = 2634.1428
$ws.Range("I77").Value = 1969.4
$ws.Range("J77").Value = 4296
$ws.Range("K77").Value = 9847
$ws.Range("L77").Value = 21480
$ws.Range("M77").Value = -5479
$ws.Range("N77").Value = -30216

# Row 88
$ws.Range("H88").Value = 1034.25
$ws.Range("J88").Value = 1069
$ws.Range("L88").Value = 1069
$ws.Range("N88").Value = -1881

# Row 91
$ws.Range("H91").Value = 1034.25
$ws.Range("J91").Value = 1069
$ws.Range("L91").Value = 1069
$ws.Range("N91").Value = -3877

# Row 110
$ws.Range("H110").Value = 2500
$ws.Range("I110").Value = 2500
$ws.Range("K110").Value = 2500
$ws.Range("M110").Value = -455

# Row 116
$ws.Range("H116").Value = 1254.4546
$ws.Range("I116").Value = 699.8889
$ws.Range("K116").Value = 699.8889
$ws.Range("M116").Value = 1594.1111

# Row 122
$ws.Range("H122").Value = 1950
$ws.Range("I122").Value = 1557.2858
$ws.Range("K122").Value = 4671.857400000001
$ws.Range("M122").Value = -2221.857400000001

# Row 132
$ws.Range("H132").Value = 2319.6
$ws.Range("I132").Value = 2200
$ws.Range("K132").Value = 6600
$ws.Range("M132").Value = -4070

# Row 133
$ws.Range("H133").Value = 500261
$ws.Range("J133").Value = 500261
$ws.Range("L133").Value = 500261
$ws.Range("N133").Value = -505321

# Row 136
$ws.Range("H136").Value = 5308.4375
$ws.Range("I136").Value = 1849.4286
$ws.Range("K136").Value = 5548.2858
$ws.Range("M136").Value = -2998.2858

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1254.4546
$ws.Range("I3").Value = 699.8889
$ws.Range("K3").Value = 699.8889
$ws.Range("M3").Value = -585.8889

# Row 22
$ws.Range("H22").Value = 709.36365
$ws.Range("I22").Value = 644.7778
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 644.7778
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -471.7778
$ws.Range("N22").Value = -1346

# Row 95
$ws.Range("H95").Value = 21171
$ws.Range("J95").Value = 21171
$ws.Range("L95").Value = 21171
$ws.Range("N95").Value = -26663

# Row 107
$ws.Range("H107").Value = 3632.6155
$ws.Range("I107").Value = 2355.3333
$ws.Range("K107").Value = 2355.3333
$ws.Range("M107").Value = -435.3332999999998

# Row 127
$ws.Range("H127").Value = 69999
$ws.Range("J127").Value = 69999
$ws.Range("L127").Value = 69999
$ws.Range("N127").Value = -79919

# Row 134
$ws.Range("H134").Value = 1806.9286
$ws.Range("I134").Value = 1806.9286
$ws.Range("K134").Value = 5420.7858
$ws.Range("M134").Value = -2885.7858

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1561.6111
$ws.Range("I16").Value = 1142.7142
$ws.Range("J16").Value = 3027.75
$ws.Range("K16").Value = 1142.7142
$ws.Range("L16").Value = 3027.75
$ws.Range("M16").Value = -855.7141999999999
$ws.Range("N16").Value = -3601.75

# Row 58
$ws.Range("H58").Value = 2067.5
$ws.Range("I58").Value = 2125.8333
$ws.Range("K58").Value = 2125.8333
$ws.Range("M58").Value = -1922.8333

# Row 94
$ws.Range("H94").Value = 94328.25
$ws.Range("I94").Value = 184686
$ws.Range("J94").Value = 3970.5
$ws.Range("K94").Value = 184686
$ws.Range("L94").Value = 3970.5
$ws.Range("M94").Value = -184235
$ws.Range("N94").Value = -4872.5

# Row 113
$ws.Range("H113").Value = 1561.6111
$ws.Range("I113").Value = 1142.7142
$ws.Range("J113").Value = 3027.75
$ws.Range("K113").Value = 1142.7142
$ws.Range("L113").Value = 3027.75
$ws.Range("M113").Value = 1027.2858
$ws.Range("N113").Value = -7367.75

# Row 134
$ws.Range("H134").Value = 3386.375
$ws.Range("I134").Value = 3584.7144
$ws.Range("J134").Value = 1998
$ws.Range("K134").Value = 10754.1432
$ws.Range("L134").Value = 5994
$ws.Range("M134").Value = -8219.143199999999
$ws.Range("N134").Value = -11064

# Row 136
$ws.Range("H136").Value = 2067.5
$ws.Range("I136").Value = 2125.8333
$ws.Range("K136").Value = 6377.499899999999
$ws.Range("M136").Value = -3827.499899999999

$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("H3").Value = 7664.6665
$ws.Range("I3").Value = 7664.6665
$ws.Range("K3").Value = 22993.9995
$ws.Range("M3").Value = -22881.9995

# Row 140
$ws.Range("H140").Value = 6837.778
$ws.Range("I140").Value = 2756.9167
$ws.Range("K140").Value = 8270.750100000001
$ws.Range("M140").Value = -3090.750100000001

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 1604.9375
$ws.Range("I102").Value = 1229.1538
$ws.Range("K102").Value = 1229.1538
$ws.Range("M102").Value = 392.8462

# Row 107
$ws.Range("H107").Value = 3912.8333
$ws.Range("J107").Value = 5435.5
$ws.Range("L107").Value = 5435.5
$ws.Range("N107").Value = -9275.5

# Row 113
$ws.Range("H113").Value = 1211.375
$ws.Range("J113").Value = 1264.6666
$ws.Range("L113").Value = 1264.6666
$ws.Range("N113").Value = -5604.6666

# Row 123
$ws.Range("H123").Value = 72500
$ws.Range("J123").Value = 72500
$ws.Range("L123").Value = 72500
$ws.Range("N123").Value = -77400

$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 2765.5715
$ws.Range("I61").Value = 2153.95
$ws.Range("K61").Value = 2153.95
$ws.Range("M61").Value = -1951.95

# Row 113
$ws.Range("H113").Value = 2765.5715
$ws.Range("I113").Value = 2153.95
$ws.Range("K113").Value = 2153.95
$ws.Range("M113").Value = 16.05000000000018

# Row 132
$ws.Range("H132").Value = 2869.8
$ws.Range("I132").Value = 1281.3334
$ws.Range("K132").Value = 3844.0002
$ws.Range("M132").Value = -1314.0002

# Row 136
$ws.Range("H136").Value = 4919.6
$ws.Range("I136").Value = 4896.5
$ws.Range("J136").Value = 4935
$ws.Range("K136").Value = 14689.5
$ws.Range("L136").Value = 14805
$ws.Range("M136").Value = -12139.5
$ws.Range("N136").Value = -19905

$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 5547.3335
$ws.Range("I126").Value = 3013.6667
$ws.Range("K126").Value = 9041.000100000001
$ws.Range("M126").Value = -6571.000100000001

# Row 132
$ws.Range("H132").Value = 1796.2
$ws.Range("I132").Value = 1796
$ws.Range("J132").Value = 1797
$ws.Range("K132").Value = 5388
$ws.Range("L132").Value = 5391
$ws.Range("M132").Value = -2858
$ws.Range("N132").Value = -10451
